$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 189, pushing existing rows 189-227 down to 190-228.
$ws.Rows.Item(189).Insert()

# Populate the newly inserted row 189 with the new record.
$ws.Range("A189").Value = 4
$ws.Range("B189").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C189").Value = "Los Lagos"
$ws.Range("D189").Value = 44551
$ws.Range("E189").Value = 10
$ws.Range("F189").Value = 100114014
$ws.Range("G189").Value = "Betarraga"
$ws.Range("H189").Value = "Sin especificar"
$ws.Range("I189").Value = "Primera"
$ws.Range("J189").Value = 1000
$ws.Range("K189").Value = 900
$ws.Range("L189").Value = 1000
$ws.Range("M189").Value = 950
$ws.Range("N189").Value = "$/paquete 5 unidades"
$ws.Range("O189").Value = "Región del Maule"
$ws.Range("P189").Value = 190
$ws.Range("Q189").Value = 5
$ws.Range("R189").Value = "Hortaliza"
